$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data area (old range was A2:B29) before writing the new table.
$ws.Range("A1:D40").ClearContents()

# Header row
$ws.Range("A2").Value = "ZIP Code"
$ws.Range("B2").Value = "Deaths_per_HundThou"
$ws.Range("C2").Value = "Cases_per_HundThou"

# Updated data table: ZIP Code, Deaths_per_HundThou, Cases_per_HundThou
$data = @(
  @(48201, 209, 1885),
  @(48202, 193, 1884),
  @(48203, 142, 1094),
  @(48204, 204, 1735),
  @(48205, 163, 1566),
  @(48206, 124, 1528),
  @(48207, 525, 2893),
  @(48208, 283, 2155),
  @(48210, 67, 1959),
  @(48211, 265, 2849),
  @(48213, 249, 1593),
  @(48214, 341, 1974),
  @(48215, 185, 1646),
  @(48219, 299, 2211),
  @(48221, 223, 2153),
  @(48223, 153, 1718),
  @(48224, 224, 1770),
  @(48227, 194, 1817),
  @(48228, 160, 1958),
  @(48234, 204, 1610),
  @(48235, 381, 2707),
  @(48238, 121, 1490)
)

$row = 3
foreach ($r in $data) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $row = $row + 1
}

$lastRow = $row - 1

# Column C width to match new sizing (closest achievable value to 18.83203125
# given this runtime's width-rounding behaviour)
$ws.Columns.Item(3).ColumnWidth = 18.0

# Selection reflects the edited row (row 11) spanning the whole row
$ws.Range("A11:XFD11").Select()

# Re-apply sort metadata (kept to the original A:B span of the sort, as in the source)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3:A" + $lastRow))
$ws.Sort.SetRange($ws.Range("A2:B" + $lastRow))
$ws.Sort.Header = 1
$ws.Sort.Apply()
